$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for each changed cell (row -> column -> value)
$updates = @{
    2  = @{ G = 0.5979736666666666;  H = 1.793921;            I = 0.03342655292740804;  J = 0.03342655292740804;
            M = 77.08952333333333;   N = 231.26857;           O = 0.2403816673726824;   P = 0.2403816673726824;
            Q = 46.09750492921889;   R = 414.87754436297;     S = 0.008035130527211564; T = 0.008035130527211564 }
    3  = @{ G = 0.5979736666666666;  H = 1.793921;            I = 0.03342655292740804;  J = 0.03342655292740804;
            O = 0.3167483425780597; P = 0.3167483425780597;
            Q = 60.74218738435022;  R = 546.679686459152;    S = 0.01058780523785429;  T = 0.01058780523785429 }
    4  = @{ G = 0.5979736666666666;  H = 1.793921;            I = 0.03342655292740804;  J = 0.03342655292740804;
            O = 0.4428699900492579; P = 0.4428699900492579;
            Q = 84.92827998254755;  R = 764.354519842928;    S = 0.01480361716234219;  T = 0.01480361716234219 }
    5  = @{ I = 0.8874158839838097; J = 0.8874158839838097;
            M = 77.08952333333333;  N = 231.26857;           O = 0.2403816673726824;   P = 0.2403816673726824;
            Q = 1223.807258111519;  R = 11014.26532300367;   S = 0.2133185098450311;   T = 0.2133185098450311 }
    6  = @{ I = 0.8874158839838097; J = 0.8874158839838097;
            O = 0.3167483425780597; P = 0.3167483425780597;
            S = 0.2810875104293155; T = 0.2810875104293154 }
    7  = @{ I = 0.8874158839838097; J = 0.8874158839838097;
            O = 0.4428699900492579; P = 0.4428699900492579;
            S = 0.3930098637094632; T = 0.3930098637094632 }
    8  = @{ I = 0.07915756308878232; J = 0.07915756308878232;
            M = 77.08952333333333;  N = 231.26857;           O = 0.2403816673726824;   P = 0.2403816673726824;
            Q = 109.1636987694934;  R = 982.4732889254401;   S = 0.01902802700043979;  T = 0.01902802700043979 }
    9  = @{ I = 0.07915756308878232; J = 0.07915756308878232;
            O = 0.3167483425780597; P = 0.3167483425780597;
            Q = 143.8438339863894;
            S = 0.02507302691089;   T = 0.02507302691088999 }
    10 = @{ I = 0.07915756308878232; J = 0.07915756308878232;
            O = 0.4428699900492579; P = 0.4428699900492579;
            Q = 201.1190234105174;
            S = 0.03505650917745253; T = 0.03505650917745253 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
